# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header fields -------------------------------------------------
# Valor Mora (total) and Cant. Trabajadores
$ws.Range("E11").Value = 779484
$ws.Range("C13").Value = 9

# --- Rebuild the worker/period detail table (rows 16-29 -> rows 16-25) ----
# First remove the 4 rows that are no longer needed (this also shifts the
# trailing signature block up from rows 34/35 to rows 30/31 automatically).
$ws.Rows("26:29").Delete()

# New data for the remaining 10 rows of the table (columns B..G).
# Column order: Tipo Doc, N Doc Trabajador, Nombre Trabajador, Periodo Mora, Valor Mora, Salario Basico
$rows = @(
    @("CC", "1010119201", "BRAYAN SAID AGRESOTT TORRENEGRA", "2505", 49066, 1600000),
    @("CC", "78712184",   "MELVIN ENRIQUE AGRESOTT PAEZ",     "2508", 160000, 4000000),
    @("CC", "8861236",    "HECTOR SANTIAGO OLIVERA ANAYA",    "2508", 56940, 1423500),
    @("CC", "1007323907", "DARWIN VEGA ROCHA",                "2508", 72250, 1806250),
    @("CC", "1047401859", "SERGIO LUIS HERNANDEZ ARAGON",     "2508", 85600, 781242),
    @("CC", "10175743",   "CESAR AUGUSTO VARGAS BUSTOS",      "2508", 67344, 1683600),
    @("CC", "72052851",   "MARCO ANTONIO CONTRERAS CAMACHO",  "2508", 67344, 1683600),
    @("CC", "1143340322", "JIMENA MARIA PUELLO PEREZ",        "2508", 100000, 2500000),
    @("CC", "1102148728", "LINDA LUCIA NARVAEZ REYES",        "2508", 56940, 1423500),
    @("CC", "1010119201", "BRAYAN SAID AGRESOTT TORRENEGRA",  "2508", 64000, 1600000)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
